$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "19.791.26"
$ws.Range("E2").Value = "  -8.48%  "
$ws.Range("D3").Value = "1.385.16"
$ws.Range("E3").Value = "  -9.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.004"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.06"
$ws.Range("E6").Value = "  -6.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3613"
$ws.Range("E7").Value = "  -8.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2998"
$ws.Range("E8").Value = "  -4.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.63"
$ws.Range("E9").Value = "  -8.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06314"
$ws.Range("E10").Value = "  -11.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9525"
$ws.Range("E11").Value = "  -8.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.007"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.220"
$ws.Range("E13").Value = "  -7.33%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.973"
$ws.Range("E14").Value = "  -9.08%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.390.45"
$ws.Range("E15").Value = "  -9.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.20"
$ws.Range("E16").Value = "  -12.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009741"
$ws.Range("E17").Value = "  -10.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05609"
$ws.Range("E18").Value = "  -14.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.52"
$ws.Range("E20").Value = "  -16.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.456"
$ws.Range("E21").Value = "  -10.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.01"
$ws.Range("E22").Value = "  -8.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.42"
$ws.Range("E23").Value = "  -3.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.252"
$ws.Range("E24").Value = "  -4.55%  "
$ws.Range("D25").Value = "19.801.65"
$ws.Range("E25").Value = "  -8.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.104"
$ws.Range("E26").Value = "  -9.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.66"
$ws.Range("E27").Value = "  -8.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.36"
$ws.Range("E28").Value = "  -10.55%  "
$ws.Range("D29").Value = "1.551.38"
$ws.Range("E29").Value = "  -9.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "107.07"
$ws.Range("E30").Value = "  -8.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.821"
$ws.Range("E31").Value = "  -21.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.165"
$ws.Range("E32").Value = "  -11.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7819"
$ws.Range("E33").Value = "  -16.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07539"
$ws.Range("E34").Value = "  -7.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.190"
$ws.Range("E35").Value = "  -3.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.004"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.624"
$ws.Range("E37").Value = "  -9.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05496"
$ws.Range("E38").Value = "  -8.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1864"
$ws.Range("E39").Value = "  -7.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01976"
$ws.Range("E40").Value = "  -9.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.293"
$ws.Range("E41").Value = "  -10.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.894"
$ws.Range("E42").Value = "  -8.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.026"
$ws.Range("E43").Value = "  -11.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.452"
$ws.Range("E44").Value = "  -7.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5111"
$ws.Range("E45").Value = "  -10.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.63"
$ws.Range("E46").Value = "  -10.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4883"
$ws.Range("E47").Value = "  -10.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.07"
$ws.Range("E48").Value = "  -7.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.693"
$ws.Range("E49").Value = "  -8.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.019"
$ws.Range("E51").Value = "  -12.40%  "
